$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the Address/modify rule text from row 6 (H6,I6) down to row 7 (H7,I7)
$ws.Range("H7").Value2 = $ws.Range("H6").Value2
$ws.Range("I7").Value2 = $ws.Range("I6").Value2
$ws.Range("H6").ClearContents()
$ws.Range("I6").ClearContents()

# Widen columns H and I (target stored widths ~41.1 / ~39.1;
# the engine quantizes ColumnWidth to pixel-based increments, so the
# input values below are chosen to land as close as possible to that)
$ws.Columns.Item(8).ColumnWidth = 40.25
$ws.Columns.Item(9).ColumnWidth = 38.25

# Update the active selection / view (topLeftCell D1, active cell I13)
$ws.Range("I13").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1 | Out-Null
